$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("study")
$ws2 = $wb.Worksheets.Item("association")

# study sheet: widen column C (MANDATORY)
$ws1.Columns.Item(3).ColumnWidth = 10.5

# association sheet: clear the autofilter criteria first (this also unhides
# the rows that used to be filtered out) before touching any cell values,
# otherwise editing still-hidden rows leaves stray row-height artifacts
$ws2.ShowAllData()

# Reorder the beta/beta_unit/ci_lower/ci_upper rows:
#   old order (12-15): beta, beta_unit, ci_lower, ci_upper
#   new order (12-15): ci_lower, ci_upper, beta, beta_unit
$ws2.Range("A12").Value = "ci_lower"
$ws2.Range("B12").Value = "Lower limit of the confidence interval"
$ws2.Range("J12").Value = "CI lower"

$ws2.Range("A13").Value = "ci_upper"
$ws2.Range("B13").Value = "Upper limit of the confidence interval"
$ws2.Range("E13").Value = "number"
$ws2.Range("J13").Value = "CI upper"

$ws2.Range("A14").Value = "beta"
$ws2.Range("B14").Value = "Effect size from a linear regression for continious traits"
$ws2.Range("J14").Value = "Beta"

$ws2.Range("A15").Value = "beta_unit"
$ws2.Range("B15").Value = "Beta unit"
$ws2.Range("E15").Value = "string"
$ws2.Range("J15").Value = "Beta unit"

# Flip MANDATORY flag (column C) to TRUE for rows 11-16
foreach ($r in 11..16) {
    $ws2.Range("C$r").Value = $true
}

# Shrink the association sheet's hidden _FilterDatabase defined name range
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "association!_FilterDatabase") {
        $dn.RefersTo = "=association!`$C`$1:`$C`$15"
    }
}

# Switch the active tab from study to association, and update the selection
$ws2.Activate() | Out-Null
$ws2.Range("A14:XFD16").Select() | Out-Null
